# Weekly data refresh for the "Hortaliza, Vega Central Mapocho de Santiago - Ají" sheet.
# A new week's observation is inserted at row 217, pushing every existing
# record (old rows 217-302) down by one row (to 218-303). The freshly
# inserted row 217 is then populated with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 217 - this shifts rows 217..302
# down to 218..303 and grows the used range from A1:R302 to A1:R303.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new observation.
$ws.Range("A217").Value = 9
$ws.Range("B217").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C217").Value = "Metropolitana"
$ws.Range("D217").Value = 44726
$ws.Range("E217").Value = 13
$ws.Range("F217").Value = 100112021
$ws.Range("G217").Value = "Ají"
$ws.Range("H217").Value = "Inferno"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 110
$ws.Range("K217").Value = 19000
$ws.Range("L217").Value = 20000
$ws.Range("M217").Value = 19364
$ws.Range("N217").Value = "$/caja 12 kilos"
$ws.Range("O217").Value = "Región de Arica y Parinacota"
$ws.Range("P217").Value = 1614
$ws.Range("Q217").Value = 12
$ws.Range("R217").Value = "Hortaliza"
